# Auto-generated Excel COM-interop script
# Applies cell-value updates to the Ixion_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 433.33334
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 433.33334
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -773.33334
$ws.Range("H21").Value = 56509.5
$ws.Range("I21").Value = 56509.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 56509.5
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -56041.5
$ws.Range("H23").Value = 56509.5
$ws.Range("I23").Value = 56509.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 56509.5
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -56275.5
$ws.Range("H28").Value = 453.18182
$ws.Range("I28").Value = 196.66667
$ws.Range("J28").Value = 1002.8571
$ws.Range("K28").Value = 196.66667
$ws.Range("L28").Value = 1002.8571
$ws.Range("M28").Value = 288.33333
$ws.Range("N28").Value = -1972.8571
$ws.Range("H112").Value = 12196070
$ws.Range("J112").Value = 14706853
$ws.Range("L112").Value = 44120559
$ws.Range("N112").Value = -44122775
$ws.Range("H113").Value = 6636.5386
$ws.Range("I113").Value = 6697.5
$ws.Range("K113").Value = 6697.5
$ws.Range("M113").Value = -3443.5
$ws.Range("H132").Value = 1786.9166
$ws.Range("I132").Value = 1691.0646
$ws.Range("J132").Value = 2381.2
$ws.Range("K132").Value = 5073.1938
$ws.Range("L132").Value = 7143.599999999999
$ws.Range("M132").Value = -2543.1938
$ws.Range("N132").Value = -12203.6
$ws.Range("H137").Value = 1745.0741
$ws.Range("I137").Value = 1656
$ws.Range("J137").Value = 1874.6364
$ws.Range("K137").Value = 4968
$ws.Range("L137").Value = 5623.9092
$ws.Range("M137").Value = -2418
$ws.Range("N137").Value = -10723.9092
$ws.Range("H138").Value = 16082.835
$ws.Range("I138").Value = 822.2593000000001
$ws.Range("J138").Value = 24006.596
$ws.Range("K138").Value = 2466.7779
$ws.Range("L138").Value = 72019.788
$ws.Range("M138").Value = 2673.2221
$ws.Range("N138").Value = -82299.788
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360
$ws.Range("H141").Value = 2316.6956
$ws.Range("I141").Value = 2225.2104
$ws.Range("J141").Value = 2751.25
$ws.Range("K141").Value = 6675.6312
$ws.Range("L141").Value = 8253.75
$ws.Range("M141").Value = -1495.6312
$ws.Range("N141").Value = -18613.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0
$ws.Range("H132").Value = 5677.921
$ws.Range("I132").Value = 1850.4348
$ws.Range("J132").Value = 11546.733
$ws.Range("K132").Value = 5551.3044
$ws.Range("L132").Value = 34640.199
$ws.Range("M132").Value = -3021.3044
$ws.Range("N132").Value = -39700.199

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 20085.5
$ws.Range("I106").Value = 19500
$ws.Range("J106").Value = 20671
$ws.Range("K106").Value = 19500
$ws.Range("L106").Value = 20671
$ws.Range("M106").Value = -18238
$ws.Range("N106").Value = -23195
$ws.Range("H107").Value = 1255.8
$ws.Range("I107").Value = 1209.8235
$ws.Range("J107").Value = 1516.3334
$ws.Range("K107").Value = 1209.8235
$ws.Range("L107").Value = 1516.3334
$ws.Range("M107").Value = 710.1765
$ws.Range("N107").Value = -5356.3334
$ws.Range("H140").Value = 48457.5
$ws.Range("J140").Value = 48457.5
$ws.Range("L140").Value = 48457.5
$ws.Range("N140").Value = -58817.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 7764.857
$ws.Range("I31").Value = 1260.4546
$ws.Range("J31").Value = 18772.309
$ws.Range("K31").Value = 1260.4546
$ws.Range("L31").Value = 18772.309
$ws.Range("M31").Value = -965.4546
$ws.Range("N31").Value = -19362.309
$ws.Range("H34").Value = 7764.857
$ws.Range("I34").Value = 1260.4546
$ws.Range("J34").Value = 18772.309
$ws.Range("K34").Value = 1260.4546
$ws.Range("L34").Value = 18772.309
$ws.Range("M34").Value = -1058.4546
$ws.Range("N34").Value = -19176.309
$ws.Range("H80").Value = 27064
$ws.Range("J80").Value = 27064
$ws.Range("L80").Value = 27064
$ws.Range("N80").Value = -29310
$ws.Range("H83").Value = 27064
$ws.Range("J83").Value = 27064
$ws.Range("L83").Value = 81192
$ws.Range("N83").Value = -92424

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7543.0234
$ws.Range("I3").Value = 14505.556
$ws.Range("J3").Value = 5700
$ws.Range("K3").Value = 43516.66800000001
$ws.Range("L3").Value = 17100
$ws.Range("M3").Value = -43404.66800000001
$ws.Range("N3").Value = -17324
$ws.Range("H5").Value = 301141
$ws.Range("I5").Value = 601.875
$ws.Range("J5").Value = 501500.4
$ws.Range("K5").Value = 1805.625
$ws.Range("L5").Value = 1504501.2
$ws.Range("M5").Value = -1693.625
$ws.Range("N5").Value = -1504725.2
$ws.Range("H122").Value = 801
$ws.Range("J122").Value = 850
$ws.Range("L122").Value = 7650
$ws.Range("N122").Value = -12550
$ws.Range("H131").Value = 35715520
$ws.Range("I131").Value = 702.93335
$ws.Range("J131").Value = 76924920
$ws.Range("K131").Value = 2108.80005
$ws.Range("L131").Value = 230774760
$ws.Range("M131").Value = 2931.19995
$ws.Range("N131").Value = -230784840
$ws.Range("H133").Value = 36110.312
$ws.Range("I133").Value = 103846.1
$ws.Range("J133").Value = 9016
$ws.Range("K133").Value = 311538.3
$ws.Range("L133").Value = 27048
$ws.Range("M133").Value = -306478.3
$ws.Range("N133").Value = -37168
$ws.Range("H135").Value = 301141
$ws.Range("I135").Value = 601.875
$ws.Range("J135").Value = 501500.4
$ws.Range("K135").Value = 5416.875
$ws.Range("L135").Value = 4513503.600000001
$ws.Range("M135").Value = -2881.875
$ws.Range("N135").Value = -4518573.600000001
$ws.Range("H136").Value = 4981.1934
$ws.Range("I136").Value = 17402.834
$ws.Range("K136").Value = 52208.50199999999
$ws.Range("M136").Value = -47108.50199999999
$ws.Range("H138").Value = 13639.091
$ws.Range("I138").Value = 17718.572
$ws.Range("J138").Value = 6500
$ws.Range("K138").Value = 53155.716
$ws.Range("L138").Value = 19500
$ws.Range("M138").Value = -48015.716
$ws.Range("N138").Value = -29780
$ws.Range("H139").Value = 4147.5557
$ws.Range("I139").Value = 5504.5454
$ws.Range("J139").Value = 2849.5652
$ws.Range("K139").Value = 16513.6362
$ws.Range("L139").Value = 8548.695599999999
$ws.Range("M139").Value = -11373.6362
$ws.Range("N139").Value = -18828.6956

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 33118
$ws.Range("J15").Value = 33118
$ws.Range("L15").Value = 33118
$ws.Range("N15").Value = -33694
$ws.Range("H81").Value = 33118
$ws.Range("J81").Value = 33118
$ws.Range("L81").Value = 33118
$ws.Range("N81").Value = -35114
$ws.Range("H84").Value = 33118
$ws.Range("J84").Value = 33118
$ws.Range("L84").Value = 99354
$ws.Range("N84").Value = -109338
$ws.Range("H102").Value = 3400.9524
$ws.Range("I102").Value = 2037.091
$ws.Range("J102").Value = 4901.2
$ws.Range("K102").Value = 2037.091
$ws.Range("L102").Value = 4901.2
$ws.Range("M102").Value = -415.0909999999999
$ws.Range("N102").Value = -8145.2
$ws.Range("H139").Value = 40326
$ws.Range("J139").Value = 40326
$ws.Range("L139").Value = 40326
$ws.Range("N139").Value = -50606

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 999
$ws.Range("I13").Value = 999
$ws.Range("K13").Value = 999
$ws.Range("M13").Value = -859
$ws.Range("H61").Value = 2660.8
$ws.Range("I61").Value = 2035.3
$ws.Range("J61").Value = 3911.8
$ws.Range("K61").Value = 2035.3
$ws.Range("L61").Value = 3911.8
$ws.Range("M61").Value = -1833.3
$ws.Range("N61").Value = -4315.8
$ws.Range("H113").Value = 2660.8
$ws.Range("I113").Value = 2035.3
$ws.Range("J113").Value = 3911.8
$ws.Range("K113").Value = 2035.3
$ws.Range("L113").Value = 3911.8
$ws.Range("M113").Value = 134.7
$ws.Range("N113").Value = -8251.799999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000000
$ws.Range("I15").Value = 10000000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 10000000
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = -9999712
$ws.Range("N15").Value = 0
$ws.Range("H122").Value = 2873.8
$ws.Range("I122").Value = 1533.8334
$ws.Range("J122").Value = 3767.111
$ws.Range("K122").Value = 4601.5002
$ws.Range("L122").Value = 11301.333
$ws.Range("M122").Value = -2151.5002
$ws.Range("N122").Value = -16201.333
$ws.Range("H136").Value = 1522
$ws.Range("I136").Value = 1173.7
$ws.Range("J136").Value = 2102.5
$ws.Range("K136").Value = 3521.1
$ws.Range("L136").Value = 6307.5
$ws.Range("M136").Value = -971.1000000000004
$ws.Range("N136").Value = -11407.5
$ws.Range("H141").Value = 60745
$ws.Range("J141").Value = 67569.164
$ws.Range("L141").Value = 67569.164
$ws.Range("N141").Value = -77929.164

